# Applies the edits described by the diff to the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows that simply lose their T/U/V (Unit_Price_USD / TOTAL_ASS_VALUE_USD /
#    Invoice_Unit_Price_FC_USD) values, with no style changes.
# ---------------------------------------------------------------------------
$simpleClearRows = @(12, 15, 16, 17, 18, 19, 31, 35, 36, 37, 38, 43)
foreach ($r in $simpleClearRows) {
    $addr = "T" + $r + ":V" + $r
    $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Row 34 keeps its (yellow-highlight) style, but T/U/V are cleared too.
# ---------------------------------------------------------------------------
$ws.Range("T34:V34").ClearContents()

# ---------------------------------------------------------------------------
# 3) Row 14: the whole row loses its special (yellow-highlight) style, the
#    Port_of_Shipment (G14) cell is cleared out entirely, and the Importer
#    name (P14) is corrected.
# ---------------------------------------------------------------------------
$ws.Range("A14:V14").Style = "Normal"
$ws.Range("G14").ClearContents()
$ws.Range("P14").Value = "ultratech cement limited"

# ---------------------------------------------------------------------------
# 4) Row 32: same kind of change as row 14 - style removed, Port_of_Shipment
#    cleared, Importer name corrected.
# ---------------------------------------------------------------------------
$ws.Range("A32:V32").Style = "Normal"
$ws.Range("G32").ClearContents()
$ws.Range("P32").Value = "brakes india private limited"

# ---------------------------------------------------------------------------
# 5) Row 46: Quantity (K46) corrected.
# ---------------------------------------------------------------------------
$ws.Range("K46").Value = 50000

# ---------------------------------------------------------------------------
# 6) Append four brand-new rows (47-50) of calcined-petroleum-coke imports.
#    (Note: this engine's PowerShell dialect does not bind named
#    "-Param value" arguments inside functions, so positional args are used.)
# ---------------------------------------------------------------------------
function Set-ImportRow(
    $ws, $Row, $Month, $Location, $BeDate, $Hscode, $Country,
    $ItemDescription, $Currency, $InvoiceUnitPriceFc, $Quantity, $Uqc,
    $UnitPrice, $TotalAssValue, $Supplier, $Importer, $ImporterCity,
    $MonthName, $Year, $UnitPriceUsd, $TotalAssValueUsd, $InvoiceUnitPriceFcUsd
) {
    # Note: column G (Port_of_Shipment) is intentionally left untouched/blank
    # for these rows, mirroring row 46 which this data was modelled on.
    $rowRangeLeft = $ws.Range("A" + $Row + ":F" + $Row)
    $rowRangeRight = $ws.Range("H" + $Row + ":V" + $Row)

    # Only the genuinely textual columns need to be forced to Text format
    # so that Excel doesn't auto-convert date-looking strings (e.g.
    # "2019-08-01") into real date serial numbers while being written.
    # The numeric columns (E,J,K,M,N,S,T,U,V) are left alone.
    $textCols = @("A", "B", "C", "D", "F", "H", "I", "L", "O", "P", "Q", "R")
    foreach ($col in $textCols) {
        $ws.Range($col + $Row).NumberFormat = "@"
    }

    $ws.Range("A" + $Row).Value = "import"
    $ws.Range("B" + $Row).Value = $Month
    $ws.Range("C" + $Row).Value = $Location
    $ws.Range("D" + $Row).Value = $BeDate
    $ws.Range("E" + $Row).Value = $Hscode
    $ws.Range("F" + $Row).Value = $Country
    $ws.Range("H" + $Row).Value = $ItemDescription
    $ws.Range("I" + $Row).Value = $Currency
    $ws.Range("J" + $Row).Value = $InvoiceUnitPriceFc
    $ws.Range("K" + $Row).Value = $Quantity
    $ws.Range("L" + $Row).Value = $Uqc
    $ws.Range("M" + $Row).Value = $UnitPrice
    $ws.Range("N" + $Row).Value = $TotalAssValue
    $ws.Range("O" + $Row).Value = $Supplier
    $ws.Range("P" + $Row).Value = $Importer
    $ws.Range("Q" + $Row).Value = $ImporterCity
    $ws.Range("R" + $Row).Value = $MonthName
    $ws.Range("S" + $Row).Value = $Year
    $ws.Range("T" + $Row).Value = $UnitPriceUsd
    $ws.Range("U" + $Row).Value = $TotalAssValueUsd
    $ws.Range("V" + $Row).Value = $InvoiceUnitPriceFcUsd

    # Restore the default (unstyled) look - the source rows this was
    # modelled on (e.g. row 46) carry no explicit style/format either.
    $rowRangeLeft.Style = "Normal"
    $rowRangeRight.Style = "Normal"
}

$newRows = @(
    @{
        Row = 47; Month = "2019-08-01"; Location = "kakinada sea (inkak1)";
        BeDate = "2018-05-25 00:00:00"; Hscode = 27131200; Country = "united states";
        ItemDescription = "calcined petroleum coke"; Currency = "usd";
        InvoiceUnitPriceFc = 485; Quantity = 100000; Uqc = "kgs";
        UnitPrice = 33296.29; TotalAssValue = 173074128.2;
        Supplier = "oxbow calcining international"; Importer = "vedanta limited";
        ImporterCity = "panajigoa"; MonthName = "january"; Year = 2021;
        UnitPriceUsd = 33296.29; TotalAssValueUsd = 173074128.2; InvoiceUnitPriceFcUsd = 485
    },
    @{
        Row = 48; Month = "2020-06-01"; Location = "kakinada sea (inkak1)";
        BeDate = "2018-05-25 00:00:00"; Hscode = 27131200; Country = "united states";
        ItemDescription = "calcined petroleum coke in bulk"; Currency = "usd";
        InvoiceUnitPriceFc = 485; Quantity = 150000; Uqc = "kgs";
        UnitPrice = 33296.29; TotalAssValue = 173074128.2;
        Supplier = "oxbow calcining international"; Importer = "vedanta limited";
        ImporterCity = "panajigoa"; MonthName = "january"; Year = 2021;
        UnitPriceUsd = 33296.29; TotalAssValueUsd = 173074128.2; InvoiceUnitPriceFcUsd = 485
    },
    @{
        Row = 49; Month = "2021-07-01"; Location = "kakinada sea (inkak1)";
        BeDate = "2018-05-25 00:00:00"; Hscode = 27131200; Country = "united states";
        ItemDescription = "calcined petroleum coke in bulk"; Currency = "usd";
        InvoiceUnitPriceFc = 485; Quantity = 180000; Uqc = "kgs";
        UnitPrice = 33296.29; TotalAssValue = 173074128.2;
        Supplier = "oxbow calcining international"; Importer = "vedanta limited";
        ImporterCity = "panajigoa"; MonthName = "january"; Year = 2021;
        UnitPriceUsd = 33296.29; TotalAssValueUsd = 173074128.2; InvoiceUnitPriceFcUsd = 485
    },
    @{
        Row = 50; Month = "2022-12-01"; Location = "kakinada sea (inkak1)";
        BeDate = "2018-05-25 00:00:00"; Hscode = 27131200; Country = "united states";
        ItemDescription = "calcined petroleum coke in bulk"; Currency = "usd";
        InvoiceUnitPriceFc = 485; Quantity = 190000; Uqc = "kgs";
        UnitPrice = 33296.29; TotalAssValue = 173074128.2;
        Supplier = "oxbow calcining international"; Importer = "vedanta limited";
        ImporterCity = "panajigoa"; MonthName = "january"; Year = 2021;
        UnitPriceUsd = 33296.29; TotalAssValueUsd = 173074128.2; InvoiceUnitPriceFcUsd = 485
    }
)

foreach ($item in $newRows) {
    Set-ImportRow $ws $item.Row $item.Month $item.Location $item.BeDate `
        $item.Hscode $item.Country $item.ItemDescription $item.Currency `
        $item.InvoiceUnitPriceFc $item.Quantity $item.Uqc $item.UnitPrice `
        $item.TotalAssValue $item.Supplier $item.Importer $item.ImporterCity `
        $item.MonthName $item.Year $item.UnitPriceUsd $item.TotalAssValueUsd `
        $item.InvoiceUnitPriceFcUsd
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
